$wb = $excel.ActiveWorkbook

# --- Sheet1 edits ---
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Cells.Item(6, 1).Value = 11111111   # A6
$ws1.Cells.Item(6, 5).Value = 333.33     # E6
$ws1.Cells.Item(6, 6).Value = 300        # F6
$ws1.Cells.Item(7, 1).Value = 24256354   # A7

# --- Sheet2 edits (only column A for rows 6-14) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Cells.Item(6, 1).Value = 22222222    # A6
$ws2.Cells.Item(7, 1).Value = 44256354    # A7
$ws2.Cells.Item(8, 1).Value = 55555555    # A8
$ws2.Cells.Item(9, 1).Value = 66666666    # A9
$ws2.Cells.Item(10, 1).Value = 77777777   # A10
$ws2.Cells.Item(11, 1).Value = 88888888   # A11
$ws2.Cells.Item(12, 1).Value = 99999999   # A12
$ws2.Cells.Item(13, 1).Value = 10101011   # A13
$ws2.Cells.Item(14, 1).Value = 11111112   # A14

# --- New Sheet3: duplicate Sheet2's layout/formatting, then fix values ---
$ws2.Copy($null, $ws2)
$ws3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3.Name = "Sheet3"

$ws3.Cells.Item(6, 1).Value = 33333333
$ws3.Cells.Item(6, 3).Value = 1
$ws3.Cells.Item(6, 5).Value = 112.17
$ws3.Cells.Item(6, 6).Value = 7.3

$ws3.Cells.Item(7, 1).Value = 55556354
$ws3.Cells.Item(7, 3).Value = 1
$ws3.Cells.Item(7, 5).Value = 200
$ws3.Cells.Item(7, 6).Value = 4.34

$ws3.Cells.Item(8, 1).Value = 65555555
$ws3.Cells.Item(8, 3).Value = 1
$ws3.Cells.Item(8, 5).Value = 114.2
$ws3.Cells.Item(8, 6).Value = 74.23

$ws3.Cells.Item(9, 1).Value = 76666666
$ws3.Cells.Item(9, 3).Value = 2
$ws3.Cells.Item(9, 5).Value = 27.67
$ws3.Cells.Item(9, 6).Value = 18.56

$ws3.Cells.Item(10, 1).Value = 87777777
$ws3.Cells.Item(10, 3).Value = 1
$ws3.Cells.Item(10, 5).Value = 505
$ws3.Cells.Item(10, 6).Value = 411.14

$ws3.Cells.Item(11, 1).Value = 98888888
$ws3.Cells.Item(11, 3).Value = 1
$ws3.Cells.Item(11, 5).Value = 800
$ws3.Cells.Item(11, 6).Value = 118.84

$ws3.Cells.Item(12, 1).Value = 10999999
$ws3.Cells.Item(12, 3).Value = 1
$ws3.Cells.Item(12, 5).Value = 32.25
$ws3.Cells.Item(12, 6).Value = 20.96

$ws3.Cells.Item(13, 1).Value = 11101011
$ws3.Cells.Item(13, 3).Value = 1
$ws3.Cells.Item(13, 5).Value = 812.01
$ws3.Cells.Item(13, 6).Value = 682.09

$ws3.Cells.Item(14, 1).Value = 12111112
$ws3.Cells.Item(14, 3).Value = 1
$ws3.Cells.Item(14, 5).Value = 28.68
$ws3.Cells.Item(14, 6).Value = 1.29

$ws1.Select()
